$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 through 56 in column G all get normalized to "Shivamogga (Shimoga)"
$ws.Range("G4:G56").Value = "Shivamogga (Shimoga)"
